$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header changes ---
$ws.Range("A1").Value = "날짜"
$ws.Range("M1").Value = "엽면적지수"
$ws.Range("N1").Value = "주간생육길이_생육상태"
$ws.Range("O1").Value = "줄기굵기_생육상태"
$ws.Range("P1").Value = "잎길이_생육상태"
$ws.Range("Q1").Value = "입폭_생육상태"
$ws.Range("R1").Value = "입수_생육상태"
$ws.Range("S1").Value = "엽면적지수_생육상태"
$ws.Range("T1").Value = "개화화방위치_생육상태"
$ws.Range("U1").Value = "꽃과줄기거리_생육상태"
$ws.Range("V1").Value = "생육상태점수"
$ws.Range("W1").Value = "생장구분"

# --- Column A: week label -> numeric date (YYYYMMDD) ---
$ws.Range("A2:A16").Value = 20180322
$ws.Range("A17:A31").Value = 20180329
$ws.Range("A32:A46").Value = 20180405
$ws.Range("A47:A61").Value = 20180411
$ws.Range("A62:A76").Value = 20180418
$ws.Range("A77:A91").Value = 20180425
$ws.Range("A92:A106").Value = 20180502
$ws.Range("A107:A121").Value = 20180509
$ws.Range("A122:A136").Value = 20180516
$ws.Range("A137:A151").Value = 20180524
$ws.Range("A152:A166").Value = 20180530
$ws.Range("A167:A181").Value = 20180606
$ws.Range("A182:A196").Value = 20180613
$ws.Range("A197:A211").Value = 20180620
$ws.Range("A212:A226").Value = 20180627
$ws.Range("A227:A241").Value = 20180705
$ws.Range("A242:A256").Value = 20180711
$ws.Range("A257:A271").Value = 20180718

# --- V/W column value corrections (growth_type_score / growth_type) ---
$ws.Range("V17").Value = 2
$ws.Range("W17").Value = 1
$ws.Range("V18").Value = 1
$ws.Range("W18").Value = 1
$ws.Range("V19").Value = 1
$ws.Range("W19").Value = 1
$ws.Range("V21").Value = 0
$ws.Range("W21").Value = 0
$ws.Range("V27").Value = -2
$ws.Range("V29").Value = 0
$ws.Range("W29").Value = 0
$ws.Range("V30").Value = -4
$ws.Range("V32").Value = 4
$ws.Range("V33").Value = 4
$ws.Range("V34").Value = -3
$ws.Range("V35").Value = 1
$ws.Range("W35").Value = 1
$ws.Range("V36").Value = 4
$ws.Range("V38").Value = 4
$ws.Range("V39").Value = 1
$ws.Range("W39").Value = 1
$ws.Range("V40").Value = 1
$ws.Range("W40").Value = 1
$ws.Range("V41").Value = -3
$ws.Range("V42").Value = 1
$ws.Range("W42").Value = 1
$ws.Range("V43").Value = -2
$ws.Range("V44").Value = 4
$ws.Range("V45").Value = 4
$ws.Range("V46").Value = 2
$ws.Range("W46").Value = 1
$ws.Range("V47").Value = 0
$ws.Range("W47").Value = 0
$ws.Range("V48").Value = 0
$ws.Range("W48").Value = 0
$ws.Range("V49").Value = 0
$ws.Range("W49").Value = 0
$ws.Range("V50").Value = -3
$ws.Range("V51").Value = -3
$ws.Range("V52").Value = 3
$ws.Range("V53").Value = 3
$ws.Range("V56").Value = 2
$ws.Range("W56").Value = 1
$ws.Range("V57").Value = 1
$ws.Range("W57").Value = 1
$ws.Range("V59").Value = -2
$ws.Range("V60").Value = 0
$ws.Range("W60").Value = 0
$ws.Range("V61").Value = 3
$ws.Range("V62").Value = 2
$ws.Range("V64").Value = 2
$ws.Range("W64").Value = 1
$ws.Range("V66").Value = -5
$ws.Range("V67").Value = 4
$ws.Range("V68").Value = 3
$ws.Range("V69").Value = 1
$ws.Range("W69").Value = 1
$ws.Range("V70").Value = 3
$ws.Range("V71").Value = 0
$ws.Range("W71").Value = 0
$ws.Range("V72").Value = 0
$ws.Range("W72").Value = 0
$ws.Range("V73").Value = -1
$ws.Range("V74").Value = 2
$ws.Range("V77").Value = 3
$ws.Range("V78").Value = 3
$ws.Range("V79").Value = -2
$ws.Range("V80").Value = -2
$ws.Range("V81").Value = 2
$ws.Range("W81").Value = 1
$ws.Range("V82").Value = -1
$ws.Range("V83").Value = 2
$ws.Range("V84").Value = -2
$ws.Range("V85").Value = 0
$ws.Range("W85").Value = 0
$ws.Range("V88").Value = -1
$ws.Range("V89").Value = -1
$ws.Range("V90").Value = 1
$ws.Range("W90").Value = 1
$ws.Range("V91").Value = 0
$ws.Range("W91").Value = 0
$ws.Range("V92").Value = 0
$ws.Range("W92").Value = 0
$ws.Range("V93").Value = 2
$ws.Range("W93").Value = 1
$ws.Range("V94").Value = -2
$ws.Range("V95").Value = 1
$ws.Range("W95").Value = 1
$ws.Range("V96").Value = 3
$ws.Range("V97").Value = 1
$ws.Range("W97").Value = 1
$ws.Range("V98").Value = 5
$ws.Range("V99").Value = -3
$ws.Range("V100").Value = 3
$ws.Range("V101").Value = -2
$ws.Range("V102").Value = -1
$ws.Range("V103").Value = 2
$ws.Range("W103").Value = 1
$ws.Range("V104").Value = 1
$ws.Range("W104").Value = 1
$ws.Range("V105").Value = 1
$ws.Range("W105").Value = 1
$ws.Range("V106").Value = -2
$ws.Range("V107").Value = -5
$ws.Range("V108").Value = 0
$ws.Range("W108").Value = 0
$ws.Range("V109").Value = -3
$ws.Range("V110").Value = -4
$ws.Range("V111").Value = -2
$ws.Range("V112").Value = -3
$ws.Range("V113").Value = 2
$ws.Range("W113").Value = 1
$ws.Range("V114").Value = -1
$ws.Range("V115").Value = 0
$ws.Range("W115").Value = 0
$ws.Range("V116").Value = 4
$ws.Range("V117").Value = -2
$ws.Range("V118").Value = 0
$ws.Range("W118").Value = 0
$ws.Range("V119").Value = -1
$ws.Range("V120").Value = 0
$ws.Range("W120").Value = 0
$ws.Range("V121").Value = -3
$ws.Range("V122").Value = -3
$ws.Range("V123").Value = -3
$ws.Range("V124").Value = -3
$ws.Range("V125").Value = -2
$ws.Range("V126").Value = -3
$ws.Range("V127").Value = -2
$ws.Range("V128").Value = 2
$ws.Range("W128").Value = 1
$ws.Range("V129").Value = -2
$ws.Range("V130").Value = -3
$ws.Range("V131").Value = 2
$ws.Range("W131").Value = 1
$ws.Range("V132").Value = -4
$ws.Range("V133").Value = -2
$ws.Range("V134").Value = -1
$ws.Range("V136").Value = -3